$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calibration")
Write-Host $ws.Name
